# Auto-generated edit script applying scheduled market-data refresh
# to the Leve profit calculation sheets (currentAveragePrice* / LevePrice* / LeveProfit*).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H54").Value = 18600
$ws.Range("I54").Value = 0
$ws.Range("J54").Value = 18600
$ws.Range("K54").Value = 0
$ws.Range("L54").Value = 18600
$ws.Range("M54").ClearContents()
$ws.Range("N54").Value = -19572

$ws.Range("H137").Value = 2877.8408
$ws.Range("I137").Value = 3697.55
$ws.Range("K137").Value = 11092.65
$ws.Range("M137").Value = -8542.650000000001

$ws.Range("H138").Value = 3425.1924
$ws.Range("I138").Value = 1489.5862
$ws.Range("J138").Value = 5865.7393
$ws.Range("K138").Value = 4468.7586
$ws.Range("L138").Value = 17597.2179
$ws.Range("M138").Value = 671.2413999999999
$ws.Range("N138").Value = -27877.2179

$ws.Range("H141").Value = 1160643.5
$ws.Range("I141").Value = 3654.25
$ws.Range("J141").Value = 3474622
$ws.Range("K141").Value = 10962.75
$ws.Range("L141").Value = 10423866
$ws.Range("M141").Value = -5782.75
$ws.Range("N141").Value = -10434226

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H26").Value = 616.2
$ws.Range("I26").Value = 616.2
$ws.Range("K26").Value = 616.2
$ws.Range("M26").Value = -286.2

$ws.Range("H61").Value = 5028.8335
$ws.Range("I61").Value = 1726.25
$ws.Range("J61").Value = 5972.4287
$ws.Range("K61").Value = 1726.25
$ws.Range("L61").Value = 5972.4287
$ws.Range("M61").Value = -1514.25
$ws.Range("N61").Value = -6396.4287

$ws.Range("H74").Value = 673.5
$ws.Range("I74").Value = 566.2778
$ws.Range("J74").Value = 811.3570999999999
$ws.Range("K74").Value = 566.2778
$ws.Range("L74").Value = 811.3570999999999
$ws.Range("M74").Value = 307.7222
$ws.Range("N74").Value = -2559.3571

$ws.Range("H77").Value = 673.5
$ws.Range("I77").Value = 566.2778
$ws.Range("J77").Value = 811.3570999999999
$ws.Range("K77").Value = 2831.389
$ws.Range("L77").Value = 4056.7855
$ws.Range("M77").Value = 1536.611
$ws.Range("N77").Value = -12792.7855

$ws.Range("H136").Value = 5028.8335
$ws.Range("I136").Value = 1726.25
$ws.Range("J136").Value = 5972.4287
$ws.Range("K136").Value = 5178.75
$ws.Range("L136").Value = 17917.2861
$ws.Range("M136").Value = -2628.75
$ws.Range("N136").Value = -23017.2861

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 1580.238
$ws.Range("J105").Value = 2072.1428
$ws.Range("L105").Value = 2072.1428
$ws.Range("N105").Value = -5566.1428

$ws.Range("H107").Value = 4840.2
$ws.Range("I107").Value = 3550.25
$ws.Range("J107").Value = 10000
$ws.Range("K107").Value = 3550.25
$ws.Range("L107").Value = 10000
$ws.Range("M107").Value = -1630.25
$ws.Range("N107").Value = -13840

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 1900.7916
$ws.Range("I99").Value = 1230.5
$ws.Range("J99").Value = 2034.85
$ws.Range("K99").Value = 1230.5
$ws.Range("L99").Value = 2034.85
$ws.Range("M99").Value = 267.5
$ws.Range("N99").Value = -5030.85

$ws.Range("H105").Value = 3462.9375
$ws.Range("I105").Value = 3150.5
$ws.Range("J105").Value = 5650
$ws.Range("K105").Value = 3150.5
$ws.Range("L105").Value = 5650
$ws.Range("M105").Value = -1403.5
$ws.Range("N105").Value = -9144

$ws.Range("H126").Value = 1900.7916
$ws.Range("I126").Value = 1230.5
$ws.Range("J126").Value = 2034.85
$ws.Range("K126").Value = 3691.5
$ws.Range("L126").Value = 6104.549999999999
$ws.Range("M126").Value = -1221.5
$ws.Range("N126").Value = -11044.55

$ws.Range("H134").Value = 3325.1875
$ws.Range("J134").Value = 4601.375
$ws.Range("L134").Value = 13804.125
$ws.Range("N134").Value = -18874.125

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H13").Value = 500
$ws.Range("I13").Value = 500
$ws.Range("J13").Value = 0
$ws.Range("K13").Value = 1500
$ws.Range("L13").Value = 0
$ws.Range("M13").Value = -1332
$ws.Range("N13").ClearContents()

$ws.Range("H50").Value = 81001.5
$ws.Range("I50").Value = 68.333336
$ws.Range("J50").Value = 103074.18
$ws.Range("K50").Value = 205.000008
$ws.Range("L50").Value = 309222.54
$ws.Range("M50").Value = 275.999992
$ws.Range("N50").Value = -310184.54

$ws.Range("H53").Value = 81001.5
$ws.Range("I53").Value = 68.333336
$ws.Range("J53").Value = 103074.18
$ws.Range("K53").Value = 205.000008
$ws.Range("L53").Value = 309222.54
$ws.Range("M53").Value = 275.999992
$ws.Range("N53").Value = -310184.54

$ws.Range("H55").Value = 2249.375
$ws.Range("I55").Value = 363.33334
$ws.Range("J55").Value = 2684.6155
$ws.Range("K55").Value = 1090.00002
$ws.Range("L55").Value = 8053.8465
$ws.Range("M55").Value = -913.0000199999999
$ws.Range("N55").Value = -8407.8465

$ws.Range("H69").Value = 10604.363
$ws.Range("I69").Value = 893
$ws.Range("J69").Value = 14246.125
$ws.Range("K69").Value = 2679
$ws.Range("L69").Value = 42738.375
$ws.Range("M69").Value = -1868
$ws.Range("N69").Value = -44360.375

$ws.Range("H72").Value = 10604.363
$ws.Range("I72").Value = 893
$ws.Range("J72").Value = 14246.125
$ws.Range("K72").Value = 8037
$ws.Range("L72").Value = 128215.125
$ws.Range("M72").Value = -3981
$ws.Range("N72").Value = -136327.125

$ws.Range("H132").Value = 2491.2632
$ws.Range("I132").Value = 1004.1539
$ws.Range("K132").Value = 9037.3851
$ws.Range("M132").Value = -6507.3851

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 4104.8823
$ws.Range("I70").Value = 3800.6155
$ws.Range("J70").Value = 5093.75
$ws.Range("K70").Value = 3800.6155
$ws.Range("L70").Value = 5093.75
$ws.Range("M70").Value = -3530.6155
$ws.Range("N70").Value = -5633.75

$ws.Range("H73").Value = 4104.8823
$ws.Range("I73").Value = 3800.6155
$ws.Range("J73").Value = 5093.75
$ws.Range("K73").Value = 3800.6155
$ws.Range("L73").Value = 5093.75
$ws.Range("M73").Value = -2864.6155
$ws.Range("N73").Value = -6965.75

$ws.Range("H102").Value = 2539.6667
$ws.Range("I102").Value = 1601.1666
$ws.Range("K102").Value = 1601.1666
$ws.Range("M102").Value = 20.83339999999998

$ws.Range("H107").Value = 832.35297
$ws.Range("I107").Value = 222.5
$ws.Range("J107").Value = 1703.5714
$ws.Range("K107").Value = 222.5
$ws.Range("L107").Value = 1703.5714
$ws.Range("M107").Value = 1697.5
$ws.Range("N107").Value = -5543.5714

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H35").Value = 5745.125
$ws.Range("I35").Value = 160.16667
$ws.Range("J35").Value = 22500
$ws.Range("K35").Value = 160.16667
$ws.Range("L35").Value = 22500
$ws.Range("M35").Value = 175.83333
$ws.Range("N35").Value = -23172

$ws.Range("H40").Value = 2966.6667
$ws.Range("I40").Value = 1000
$ws.Range("J40").Value = 3360
$ws.Range("K40").Value = 1000
$ws.Range("L40").Value = 3360
$ws.Range("M40").Value = -864
$ws.Range("N40").Value = -3632

$ws.Range("H106").Value = 15592.5
$ws.Range("J106").Value = 15592.5
$ws.Range("L106").Value = 15592.5
$ws.Range("N106").Value = -18116.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H70").Value = 90105
$ws.Range("J70").Value = 90105
$ws.Range("L70").Value = 90105
$ws.Range("N70").Value = -90735

$ws.Range("H73").Value = 90105
$ws.Range("J73").Value = 90105
$ws.Range("L73").Value = 90105
$ws.Range("N73").Value = -92289

$ws.Range("H122").Value = 3070
$ws.Range("I122").Value = 2236.3635
$ws.Range("J122").Value = 4598.3335
$ws.Range("K122").Value = 6709.0905
$ws.Range("L122").Value = 13795.0005
$ws.Range("M122").Value = -4259.0905
$ws.Range("N122").Value = -18695.0005

$ws.Range("H126").Value = 1954.6061
$ws.Range("I126").Value = 1431.909
$ws.Range("J126").Value = 3000
$ws.Range("K126").Value = 4295.727000000001
$ws.Range("L126").Value = 9000
$ws.Range("M126").Value = -1825.727000000001
$ws.Range("N126").Value = -13940

$ws.Range("H132").Value = 7071.8213
$ws.Range("I132").Value = 2800.6
$ws.Range("K132").Value = 8401.799999999999
$ws.Range("M132").Value = -5871.799999999999
